$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update country ranking data and "last updated" timestamp for Pais sheet
# (values reconstructed from the target diff: re-ranked countries plus refreshed case counts)

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 23 de Marzo de 2020 a las 16:16'
$ws.Cells.Item(6, 2).Value = 35230
$ws.Cells.Item(6, 3).Value = 1684
$ws.Cells.Item(6, 4).Value = 184
$ws.Cells.Item(6, 5).Value = 34587
$ws.Cells.Item(15, 2).Value = 4018
$ws.Cells.Item(15, 3).Value = 436
$ws.Cells.Item(15, 5).Value = 3988
$ws.Cells.Item(17, 2).Value = 2547
$ws.Cells.Item(17, 3).Value = 162
$ws.Cells.Item(17, 5).Value = 2531
$ws.Cells.Item(17, 7).Value = 3
$ws.Cells.Item(17, 8).Value = 10
$ws.Cells.Item(22, 1).Value = 'Canada'
$ws.Cells.Item(22, 2).Value = 1550
$ws.Cells.Item(22, 3).Value = 80
$ws.Cells.Item(22, 4).Value = 18
$ws.Cells.Item(22, 5).Value = 1511
$ws.Cells.Item(22, 6).Value = 1
$ws.Cells.Item(22, 7).Value = 1
$ws.Cells.Item(22, 8).Value = 21
$ws.Cells.Item(24, 1).Value = 'Dinamarca'
$ws.Cells.Item(24, 2).Value = 1450
$ws.Cells.Item(24, 3).Value = 55
$ws.Cells.Item(24, 4).Value = 1
$ws.Cells.Item(24, 5).Value = 1425
$ws.Cells.Item(24, 6).Value = 55
$ws.Cells.Item(24, 7).Value = 11
$ws.Cells.Item(24, 8).Value = 24
$ws.Cells.Item(33, 4).Value = 11
$ws.Cells.Item(33, 5).Value = 733
$ws.Cells.Item(45, 1).Value = 'India'
$ws.Cells.Item(45, 2).Value = 467
$ws.Cells.Item(45, 3).Value = 71
$ws.Cells.Item(45, 4).Value = 24
$ws.Cells.Item(45, 5).Value = 435
$ws.Cells.Item(45, 6).Value = 0
$ws.Cells.Item(45, 7).Value = 1
$ws.Cells.Item(45, 8).Value = 8
$ws.Cells.Item(46, 1).Value = 'Filipinas'
$ws.Cells.Item(46, 2).Value = 462
$ws.Cells.Item(46, 3).Value = 82
$ws.Cells.Item(46, 4).Value = 18
$ws.Cells.Item(46, 5).Value = 411
$ws.Cells.Item(46, 6).Value = 1
$ws.Cells.Item(46, 7).Value = 8
$ws.Cells.Item(46, 8).Value = 33
$ws.Cells.Item(47, 1).Value = 'Eslovenia'
$ws.Cells.Item(47, 2).Value = 442
$ws.Cells.Item(47, 3).Value = 28
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(47, 5).Value = 439
$ws.Cells.Item(47, 6).Value = 12
$ws.Cells.Item(47, 7).Value = 1
$ws.Cells.Item(47, 8).Value = 3
$ws.Cells.Item(48, 1).Value = 'Rusia'
$ws.Cells.Item(48, 2).Value = 438
$ws.Cells.Item(48, 3).Value = 71
$ws.Cells.Item(48, 4).Value = 17
$ws.Cells.Item(48, 5).Value = 420
$ws.Cells.Item(48, 7).Value = 0
$ws.Cells.Item(48, 8).Value = 1
$ws.Cells.Item(118, 1).Value = 'Honduras'
$ws.Cells.Item(118, 3).Value = 1
$ws.Cells.Item(119, 1).Value = 'Bolivia'
$ws.Cells.Item(119, 3).Value = 3
$ws.Cells.Item(124, 1).Value = 'Paraguay'
$ws.Cells.Item(124, 3).Value = 0
$ws.Cells.Item(124, 6).Value = 1
$ws.Cells.Item(124, 7).Value = 0
$ws.Cells.Item(125, 1).Value = 'Montenegro'
$ws.Cells.Item(125, 3).Value = 1
$ws.Cells.Item(125, 6).Value = 0
$ws.Cells.Item(125, 7).Value = 1
$ws.Cells.Item(135, 1).Value = 'Islas Virgenes de los Estados Unidos'
$ws.Cells.Item(135, 2).Value = 17
$ws.Cells.Item(135, 3).Value = 11
$ws.Cells.Item(135, 5).Value = 17
$ws.Cells.Item(136, 1).Value = 'Kenia'
$ws.Cells.Item(136, 3).Value = 1
$ws.Cells.Item(137, 1).Value = 'Kirguistan'
$ws.Cells.Item(137, 2).Value = 16
$ws.Cells.Item(137, 3).Value = 2
$ws.Cells.Item(137, 4).Value = 0
$ws.Cells.Item(137, 5).Value = 16
$ws.Cells.Item(138, 1).Value = 'Gibraltar'
$ws.Cells.Item(138, 2).Value = 15
$ws.Cells.Item(138, 5).Value = 10
$ws.Cells.Item(139, 1).Value = 'Maldivas'
$ws.Cells.Item(139, 2).Value = 13
$ws.Cells.Item(139, 4).Value = 5
$ws.Cells.Item(139, 5).Value = 8
$ws.Cells.Item(140, 1).Value = 'Tanzania'
$ws.Cells.Item(140, 3).Value = 0
$ws.Cells.Item(141, 1).Value = 'Madagascar'
$ws.Cells.Item(141, 2).Value = 12
$ws.Cells.Item(141, 3).Value = 9
$ws.Cells.Item(141, 5).Value = 12
$ws.Cells.Item(142, 1).Value = 'Etiopia'
$ws.Cells.Item(142, 2).Value = 11
$ws.Cells.Item(142, 5).Value = 11
$ws.Cells.Item(143, 1).Value = 'Mongolia'
$ws.Cells.Item(143, 2).Value = 10
$ws.Cells.Item(143, 3).Value = 0
$ws.Cells.Item(143, 5).Value = 10
$ws.Cells.Item(144, 1).Value = 'Guinea Ecuatorial'
$ws.Cells.Item(144, 3).Value = 3
$ws.Cells.Item(144, 4).Value = 0
$ws.Cells.Item(144, 5).Value = 9
$ws.Cells.Item(145, 1).Value = 'Aruba'
$ws.Cells.Item(145, 2).Value = 9
$ws.Cells.Item(145, 3).Value = 0
$ws.Cells.Item(145, 4).Value = 1
$ws.Cells.Item(146, 1).Value = 'Nueva Caledonia'
$ws.Cells.Item(146, 2).Value = 8
$ws.Cells.Item(146, 3).Value = 4
$ws.Cells.Item(146, 5).Value = 8
$ws.Cells.Item(147, 1).Value = 'Seychelles'
$ws.Cells.Item(147, 2).Value = 7
$ws.Cells.Item(147, 5).Value = 7
$ws.Cells.Item(148, 1).Value = 'Bermudas'
$ws.Cells.Item(150, 1).Value = 'Isla de Man'
$ws.Cells.Item(151, 1).Value = 'Surinam'
$ws.Cells.Item(152, 1).Value = 'San Martin (Parte Francesa)'
$ws.Cells.Item(154, 1).Value = 'Bahamas'
$ws.Cells.Item(154, 3).Value = 0
$ws.Cells.Item(155, 1).Value = 'Groenlandia'
$ws.Cells.Item(157, 1).Value = 'Guinea'
$ws.Cells.Item(157, 3).Value = 2
$ws.Cells.Item(159, 1).Value = 'Zambia'
$ws.Cells.Item(160, 1).Value = 'Fiyi'
$ws.Cells.Item(160, 3).Value = 1
$ws.Cells.Item(161, 1).Value = 'San Bartolome'
$ws.Cells.Item(162, 1).Value = 'El Salvador'
$ws.Cells.Item(162, 3).Value = 0
$ws.Cells.Item(163, 1).Value = 'Republica de Africa Central'
$ws.Cells.Item(164, 1).Value = 'Liberia'
$ws.Cells.Item(165, 1).Value = 'Namibia'
$ws.Cells.Item(166, 1).Value = 'Congo'
$ws.Cells.Item(174, 1).Value = 'Butan'
$ws.Cells.Item(175, 1).Value = 'Nicaragua'
$ws.Cells.Item(176, 1).Value = 'Santa Lucia'
$ws.Cells.Item(177, 1).Value = 'Gambia'
$ws.Cells.Item(177, 4).Value = 0
$ws.Cells.Item(177, 7).Value = 1
$ws.Cells.Item(177, 8).Value = 1
$ws.Cells.Item(178, 1).Value = 'Sudan'
$ws.Cells.Item(178, 3).Value = 0
$ws.Cells.Item(178, 7).Value = 0
$ws.Cells.Item(179, 1).Value = 'Nepal'
$ws.Cells.Item(179, 3).Value = 1
$ws.Cells.Item(179, 4).Value = 1
$ws.Cells.Item(179, 8).Value = 0
$ws.Cells.Item(180, 1).Value = 'Siria'
$ws.Cells.Item(181, 1).Value = 'San Vicente y las Granadinas'
$ws.Cells.Item(182, 1).Value = 'Eritrea'
$ws.Cells.Item(183, 1).Value = 'Montserrat'
$ws.Cells.Item(183, 3).Value = 0
$ws.Cells.Item(184, 1).Value = 'Granada'
$ws.Cells.Item(185, 1).Value = 'Republica del Chad'
$ws.Cells.Item(186, 1).Value = 'Timor Oriental'
$ws.Cells.Item(187, 1).Value = 'Papua Nueva Guinea'
$ws.Cells.Item(188, 1).Value = 'Santa Sede'
$ws.Cells.Item(189, 1).Value = 'Somalia'
$ws.Cells.Item(190, 1).Value = 'San Martin (Parte Holandesa)'
$ws.Cells.Item(191, 1).Value = 'Republica de Yibuti'
$ws.Cells.Item(192, 1).Value = 'Antigua y Barbuda'
$ws.Cells.Item(193, 1).Value = 'Mozambique'
$ws.Cells.Item(194, 1).Value = 'Dominica'
$ws.Cells.Item(195, 1).Value = 'Uganda'
$ws.Cells.Item(196, 1).Value = 'Islas Turcas y Caicos'
$ws.Cells.Item(196, 3).Value = 1
